$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Sunrise Breakfast Helper (row 2) is no longer required on Sunday.
$ws.Range("E2").Value = $false

# The "Workcrew" role category was renamed to "Non-program" throughout
# the Require column (column D).
$dRange = $ws.Range("D1:D62")
foreach ($cell in $dRange.Cells) {
    $val = $cell.Value2
    if ($val -ne $null -and $val -like "*Workcrew*") {
        $cell.Value = $val -replace "Workcrew", "Non-program"
    }
}
